# Update Release-Notes.xlsx per folder-inventory refresh
$wb = $excel.ActiveWorkbook

# ---- Folder Inventory sheet ----
$inv = $wb.Worksheets.Item("Folder Inventory")

# The folder "Build Intelligent Apps with Microsoft's Copilot stack & Azure OpenAI"
# was touched again, so it now becomes the most-recently-updated entry.
# Insert a fresh row at the top of the data (row 2) and push everything else down.
$inv.Rows.Item(2).Insert()

$newTitle = "Build Intelligent Apps with Microsoft's Copilot stack & Azure OpenAI  "
$inv.Cells.Item(2, 1).Value = $newTitle
$inv.Cells.Item(2, 2).Value = $newTitle
$inv.Cells.Item(2, 3).Value = "2025-06-16 16:14:06 +0530"
$inv.Cells.Item(2, 4).Value = 1
$inv.Cells.Item(2, 5).Value = "Root"

# The same folder's previous entry (now pushed down to row 62) is obsolete - remove it.
$inv.Rows.Item(62).Delete()

# ---- Metadata sheet ----
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(3, 2).Value = "2025-06-16 10:56:17 UTC"
# Leading apostrophe keeps the numeric-looking value stored as text (as it was before).
$meta.Cells.Item(5, 2).Value = "'26"

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(5, 2).Value = "2025-06-16 16:14:06 +0530"
